$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update last_edited_time (column D) for the rows that were touched
$ws.Range("D4").Value = "2024-07-25T15:02:00.000Z"
$ws.Range("D5").Value = "2024-07-25T15:02:00.000Z"
$ws.Range("D6").Value = "2024-07-25T15:02:00.000Z"
$ws.Range("D8").Value = "2024-07-25T15:02:00.000Z"
$ws.Range("D12").Value = "2024-07-25T15:02:00.000Z"
$ws.Range("D13").Value = "2024-07-25T15:02:00.000Z"

# Update the "Tháng 7" totals on row 13 (remove the discount-rate summary offsets)
$ws.Range("AA13").Value = 207418000
$ws.Range("AE13").Value = 295430000
$ws.Range("AH13").Value = 246730000
$ws.Range("AK13").Value = 38
$ws.Range("AN13").Value = 48700000
$ws.Range("AQ13").Value = 282530000
